$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: C11 placement corrections
$ws.Range("B4").Value = 97.27
$ws.Range("C4").Value = 62.44

# Row 26: C37 placement corrections
$ws.Range("B26").Value = 94.53
$ws.Range("C26").Value = 62.51

# Rows 74-84: corrected rotation values (shifted by one row)
$ws.Range("E74").Value = 180
$ws.Range("E75").Value = 0
$ws.Range("E76").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("E78").Value = 180
$ws.Range("E79").Value = 270
$ws.Range("E80").Value = 90
$ws.Range("E81").Value = 90
$ws.Range("E82").Value = 180
$ws.Range("E83").Value = 0
$ws.Range("E84").Value = 180

# Remove the spurious S2 row (row 154) entirely
$ws.Rows(154).Delete()

# Update selection / view to match author's saved state
$excel.ActiveWindow.ScrollRow = 64
$ws.Range("G81").Select()
